$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: drop the "Doing" (C) category value, keep height at 60 ---
$ws.Range("C2").Clear()

# --- Row 3: category text moves from C to E, C gets a new value, row grows taller ---
$ws.Range("C3").Value = "نظرات محصولات و بلاگ و انجمن ادمین"
$ws.Range("E3").Value = "مدل صفحه ی اصلی کامل"
$ws.Rows(3).RowHeight = 75

# --- Row 4: drop old C value, add new E value ---
$ws.Range("C4").Clear()
$ws.Range("E4").Value = "سئو کامل محصول و بلاگ"

# --- Row 5: add new E value, row grows taller ---
$ws.Range("E5").Value = "بهینه سازی جستجو برای گروه های دارای محصول "
$ws.Rows(5).RowHeight = 90

# --- Row 6: add new E value, row grows taller ---
$ws.Range("E6").Value = "نام گروه در لیست گروه های محصول در صفحه افزودن محصول"
$ws.Rows(6).RowHeight = 105

# --- Row 7: add new E value (row height unchanged) ---
$ws.Range("E7").Value = "وضعیت موجود یا ناموجود در لیست نمایش محصولات"

# --- Row 16: drop B value, row shrinks ---
$ws.Range("B16").Clear()
$ws.Rows(16).RowHeight = 30

# --- Row 18: drop B value, row shrinks ---
$ws.Range("B18").Clear()
$ws.Rows(18).RowHeight = 45

# --- Row 19: drop B value, row shrinks ---
$ws.Range("B19").Clear()
$ws.Rows(19).RowHeight = 45

# --- Move the active selection to B3 (home/base index) ---
$null = $ws.Range("B3").Select()
